$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "261.17"
Set-TextCell "E2" "2.03%"
Set-TextCell "D3" "27.30"
Set-TextCell "E3" "1.47%"
Set-TextCell "D4" "4.711"
Set-TextCell "E4" "8.75%"
Set-TextCell "E5" "3.32%"
Set-TextCell "D6" "6.665"
Set-TextCell "E6" "0.73%"
Set-TextCell "D7" "0.8474"
Set-TextCell "E7" "-0.32%"
Set-TextCell "D8" "0.9170"
Set-TextCell "E8" "-1.57%"
Set-TextCell "D9" "0.1406"
Set-TextCell "E9" "1.98%"
Set-TextCell "D10" "0.04830"
Set-TextCell "E10" "5.56%"
Set-TextCell "D11" "0.07090"
Set-TextCell "E11" "0.10%"
Set-TextCell "D12" "0.03147"
Set-TextCell "E12" "2.60%"
Set-TextCell "E13" "-0.10%"
Set-TextCell "D14" "0.001532"
Set-TextCell "E14" "0.73%"
Set-TextCell "D15" "0.0006078"
Set-TextCell "E15" "0.28%"
Set-TextCell "D16" "0.006135"
Set-TextCell "E16" "-0.41%"
Set-TextCell "E17" "-1.03%"
Set-TextCell "D18" "3.148"
Set-TextCell "E18" "-0.61%"
Set-TextCell "D19" "2.176"
Set-TextCell "E19" "-1.28%"
Set-TextCell "E21" "2.13%"
Set-TextCell "D22" "4.099"
Set-TextCell "E22" "4.42%"
Set-TextCell "D23" "0.04262"
Set-TextCell "E23" "0.25%"
Set-TextCell "D24" "0.001220"
Set-TextCell "E24" "-0.14%"
Set-TextCell "E25" "-8.69%"
Set-TextCell "E26" "0.00%"
Set-TextCell "E27" "3.09%"
Set-TextCell "D40" "0.03875"
Set-TextCell "E40" "1.84%"
Set-TextCell "D41" "0.1113"
Set-TextCell "E41" "1.32%"
Set-TextCell "E42" "-33.87%"
Set-TextCell "E43" "21.22%"
Set-TextCell "E44" "0.36%"
Set-TextCell "D45" "0.00005325"
Set-TextCell "E45" "-0.73%"
Set-TextCell "E47" "-0.81%"
Set-TextCell "D48" "0.1353"
Set-TextCell "E48" "-46.24%"
